$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data rows for the two removed faculty entries (rows 3 and 4),
# but keep the D column cells (with their hyperlink style) present and empty.
$ws.Range("A3:C4").ClearContents()
$ws.Range("D3:D4").Hyperlinks.Delete()
$ws.Range("D3:D4").ClearContents()

# Update the remaining faculty member's email hyperlink/address and display text.
$ws.Range("D2").Value = "17it051@charusat.edu.in"
$ws.Hyperlinks.Item(1).Address = "mailto:17it051@charusat.edu.in"

# Update the active selection as recorded in the sheet view.
$ws.Range("G11").Select()
